$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 6) with the "ruhsat" (license) scenario data.
$ws.Range("A6").Value = "ruhsat"
$ws.Range("B6").Value = "Ruhsat Hatası"
$ws.Range("C6").Value = "Görselde ruhsat örneği yer almaktadır."
$ws.Range("D6").Value = "Görseldeki belge ile karşılaştırın."
$ws.Range("E6").Value = "Destek"
$ws.Range("F6").Value = "Ekran Alıntısı.jpg"
